# Update the "Tipo 1" and "Tipo 2" blocks on Sheet1 with new counts
# (commit: "funzionante BBe BR -> testare il resto")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Tipo 1 block (rows 2-13) ---
$ws.Range("B2").Value = 3    # n tipo
$ws.Range("L2").Value = 5
$ws.Range("B3").Value = 0    # n HVAC per EH
$ws.Range("B4").Value = 0    # n BB per EH
$ws.Range("B5").Value = 0    # n BR per BB
$ws.Range("B6").Value = 0    # n C-BESSHD per PH
$ws.Range("B8").Value = 0    # EH05HD
$ws.Range("B9").Value = 0    # PH2HD

# --- Tipo 2 block (rows 15-26) ---
$ws.Range("B15").Value = 2   # n tipo
$ws.Range("B16").Value = 10  # n HVAC per EH
$ws.Range("B17").Value = 2   # n BB per EH
$ws.Range("B18").Value = 10  # n BR per BB
$ws.Range("B19").Value = 2   # n C-BESSHD per PH
$ws.Range("B21").Value = 2   # EH05HD
$ws.Range("B22").Value = 1   # PH2HD

# Leave the cursor where the author's session ended up
$ws.Range("B23:I23").Select()
